$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Unprotect()
}

$wb.Worksheets.Item('!!Compartment').Range('A1').Value = '!!!ObjTables objTablesVersion=''0.0.8'' date=''2020-03-09 15:32:03'''
$wb.Worksheets.Item('!!Compartment').Range('A2').Value = '!!ObjTables type=''Data'' id=''Compartment'' name=''Compartment'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Compound').Range('A1').Value = '!!ObjTables type=''Data'' id=''Compound'' name=''Compound'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Definition').Range('A1').Value = '!!ObjTables type=''Data'' id=''Definition'' name=''Definition'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Enzyme').Range('A1').Value = '!!ObjTables type=''Data'' id=''Enzyme'' name=''Enzyme'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!FbcObjective').Range('A1').Value = '!!ObjTables type=''Data'' id=''FbcObjective'' name=''FbcObjective'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Gene').Range('A1').Value = '!!ObjTables type=''Data'' id=''Gene'' name=''Gene'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Layout').Range('A1').Value = '!!ObjTables type=''Data'' id=''Layout'' name=''Layout'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Measurement').Range('A1').Value = '!!ObjTables type=''Data'' id=''Measurement'' name=''Measurement'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!PbConfig').Range('A1').Value = '!!ObjTables type=''Data'' id=''PbConfig'' name=''PbConfig'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Position').Range('A1').Value = '!!ObjTables type=''Data'' id=''Position'' name=''Position'' date=''2020-03-09 15:32:03'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Protein').Range('A1').Value = '!!ObjTables type=''Data'' id=''Protein'' name=''Protein'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Quantity').Range('A1').Value = '!!ObjTables type=''Data'' id=''Quantity'' name=''Quantity'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' level=''1.0'' tableFormat=''row'' version=''0.1'''
$wb.Worksheets.Item('!!QuantityInfo').Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityInfo'' name=''QuantityInfo'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!QuantityMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''QuantityMatrix'' name=''QuantityMatrix'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Reaction').Range('A1').Value = '!!ObjTables type=''Data'' id=''Reaction'' name=''Reaction'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!ReactionStoichiometry').Range('A1').Value = '!!ObjTables type=''Data'' id=''ReactionStoichiometry'' name=''ReactionStoichiometry'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Regulator').Range('A1').Value = '!!ObjTables type=''Data'' id=''Regulator'' name=''Regulator'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Relation').Range('A1').Value = '!!ObjTables type=''Data'' id=''Relation'' name=''Relation'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!Relationship').Range('A1').Value = '!!ObjTables type=''Data'' id=''Relationship'' name=''Relationship'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrix'' name=''SparseMatrix'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixColumn').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixColumn'' name=''SparseMatrixColumn'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixOrdered').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixOrdered'' name=''SparseMatrixOrdered'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!SparseMatrixRow').Range('A1').Value = '!!ObjTables type=''Data'' id=''SparseMatrixRow'' name=''SparseMatrixRow'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!StoichiometricMatrix').Range('A1').Value = '!!ObjTables type=''Data'' id=''StoichiometricMatrix'' name=''StoichiometricMatrix'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!rxnconContingencyList').Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconContingencyList'' name=''rxnconContingencyList'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''
$wb.Worksheets.Item('!!rxnconReactionList').Range('A1').Value = '!!ObjTables type=''Data'' id=''rxnconReactionList'' name=''rxnconReactionList'' date=''2020-03-09 15:32:04'' objTablesVersion=''0.0.8'' tableFormat=''row'''

foreach ($ws in $wb.Worksheets) {
    $ws.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $true, $false, $false, $true, $false, $false, $false)
}
